$wb = $excel.ActiveWorkbook

# ---------- ALC ----------
$ws = $wb.Worksheets.Item("ALC")

# Row 93
$ws.Range("H93").Value = 21052
$ws.Range("I93").Value = 21052
$ws.Range("K93").Value = 21052
$ws.Range("M93").Value = -18556

# Row 137
$ws.Range("H137").Value = 3370.5688
$ws.Range("I137").Value = 1314.742
$ws.Range("J137").Value = 5730.963
$ws.Range("K137").Value = 3944.226
$ws.Range("L137").Value = 17192.889
$ws.Range("M137").Value = -1394.226
$ws.Range("N137").Value = -22292.889

# Row 141
$ws.Range("H141").Value = 6938.4375
$ws.Range("I141").Value = 3872.2917
$ws.Range("J141").Value = 16136.875
$ws.Range("K141").Value = 11616.8751
$ws.Range("L141").Value = 48410.625
$ws.Range("M141").Value = -6436.8751
$ws.Range("N141").Value = -58770.625

# ---------- ARM ----------
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 247490.7
$ws.Range("I32").Value = 258063
$ws.Range("J32").Value = 14900
$ws.Range("K32").Value = 258063
$ws.Range("L32").Value = 14900
$ws.Range("M32").Value = -257776
$ws.Range("N32").Value = -15474

# Row 74
$ws.Range("H74").Value = 3484.6445
$ws.Range("I74").Value = 872
$ws.Range("J74").Value = 11560.091
$ws.Range("K74").Value = 872
$ws.Range("L74").Value = 11560.091
$ws.Range("M74").Value = 2
$ws.Range("N74").Value = -13308.091

# Row 77
$ws.Range("H77").Value = 3484.6445
$ws.Range("I77").Value = 872
$ws.Range("J77").Value = 11560.091
$ws.Range("K77").Value = 4360
$ws.Range("L77").Value = 57800.455
$ws.Range("M77").Value = 8
$ws.Range("N77").Value = -66536.455

# Row 82 (M82 newly added)
$ws.Range("H82").Value = 45666.668
$ws.Range("I82").Value = 45000
$ws.Range("K82").Value = 45000
$ws.Range("M82").Value = -44639

# Row 85 (M85 newly added)
$ws.Range("H85").Value = 45666.668
$ws.Range("I85").Value = 45000
$ws.Range("K85").Value = 45000
$ws.Range("M85").Value = -43752

# Row 92
$ws.Range("H92").Value = 26509.6
$ws.Range("J92").Value = 26509.6
$ws.Range("L92").Value = 26509.6
$ws.Range("N92").Value = -31501.6

# Row 105 (N105 newly added)
$ws.Range("H105").Value = 47500
$ws.Range("J105").Value = 47500
$ws.Range("L105").Value = 47500
$ws.Range("N105").Value = -54488

# ---------- BSM ----------
$ws = $wb.Worksheets.Item("BSM")

# Row 92
$ws.Range("H92").Value = 17975
$ws.Range("J92").Value = 17975
$ws.Range("L92").Value = 17975
$ws.Range("N92").Value = -22967

# ---------- CRP ----------
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 1867.07
$ws.Range("I31").Value = 1159.8511
$ws.Range("J31").Value = 2494.2263
$ws.Range("K31").Value = 1159.8511
$ws.Range("L31").Value = 2494.2263
$ws.Range("M31").Value = -864.8511
$ws.Range("N31").Value = -3084.2263

# Row 33 (N33 removed)
$ws.Range("H33").Value = 1746.0714
$ws.Range("I33").Value = 1746.0714
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1746.0714
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1367.0714
$ws.Range("N33").Value = ""

# Row 34
$ws.Range("H34").Value = 1867.07
$ws.Range("I34").Value = 1159.8511
$ws.Range("J34").Value = 2494.2263
$ws.Range("K34").Value = 1159.8511
$ws.Range("L34").Value = 2494.2263
$ws.Range("M34").Value = -957.8511
$ws.Range("N34").Value = -2898.2263

# Row 43
$ws.Range("H43").Value = 23542.857
$ws.Range("J43").Value = 23542.857
$ws.Range("L43").Value = 23542.857
$ws.Range("N43").Value = -23910.857

# Row 44 (M44 removed)
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").Value = ""

# Row 80
$ws.Range("H80").Value = 25666.666
$ws.Range("J80").Value = 25666.666
$ws.Range("L80").Value = 25666.666
$ws.Range("N80").Value = -27912.666

# Row 83
$ws.Range("H83").Value = 25666.666
$ws.Range("J83").Value = 25666.666
$ws.Range("L83").Value = 76999.998
$ws.Range("N83").Value = -88231.998

# Row 92
$ws.Range("H92").Value = 18250
$ws.Range("J92").Value = 18250
$ws.Range("L92").Value = 18250
$ws.Range("N92").Value = -23242

# Row 95
$ws.Range("H95").Value = 17070.666
$ws.Range("J95").Value = 17070.666
$ws.Range("L95").Value = 17070.666
$ws.Range("N95").Value = -22562.666

# Row 101
$ws.Range("H101").Value = 23542.857
$ws.Range("J101").Value = 23542.857
$ws.Range("L101").Value = 23542.857
$ws.Range("N101").Value = -30032.857

# ---------- CUL ----------
$ws = $wb.Worksheets.Item("CUL")

# Row 131
$ws.Range("H131").Value = 693.7692
$ws.Range("I131").Value = 346.66666
$ws.Range("J131").Value = 991.2857
$ws.Range("K131").Value = 1039.99998
$ws.Range("L131").Value = 2973.8571
$ws.Range("M131").Value = 4000.00002
$ws.Range("N131").Value = -13053.8571

# ---------- GSM ----------
$ws = $wb.Worksheets.Item("GSM")

# Row 15
$ws.Range("H15").Value = 20464
$ws.Range("J15").Value = 20464
$ws.Range("L15").Value = 20464
$ws.Range("N15").Value = -21040

# Row 81
$ws.Range("H81").Value = 20464
$ws.Range("J81").Value = 20464
$ws.Range("L81").Value = 20464
$ws.Range("N81").Value = -22460

# Row 84
$ws.Range("H84").Value = 20464
$ws.Range("J84").Value = 20464
$ws.Range("L84").Value = 61392
$ws.Range("N84").Value = -71376

# Row 95
$ws.Range("H95").Value = 21172.25
$ws.Range("J95").Value = 21172.25
$ws.Range("L95").Value = 21172.25
$ws.Range("N95").Value = -26664.25

# Row 101
$ws.Range("H101").Value = 29300
$ws.Range("J101").Value = 29300
$ws.Range("L101").Value = 29300
$ws.Range("N101").Value = -35790

# ---------- LTW ----------
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 6098.8887
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 6798.75
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 6798.75
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -7388.75

# Row 27
$ws.Range("H27").Value = 6098.8887
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 6798.75
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 6798.75
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -7012.75

# ---------- WVR ----------
$ws = $wb.Worksheets.Item("WVR")

# Row 12
$ws.Range("H12").Value = 3068.6667
$ws.Range("I12").Value = 2956
$ws.Range("J12").Value = 3125
$ws.Range("K12").Value = 2956
$ws.Range("L12").Value = 3125
$ws.Range("M12").Value = -2814
$ws.Range("N12").Value = -3409

# Row 93
$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -44992

# Row 95
$ws.Range("H95").Value = 18277.666
$ws.Range("J95").Value = 18277.666
$ws.Range("L95").Value = 18277.666
$ws.Range("N95").Value = -23769.666

# Row 103
$ws.Range("H103").Value = 19806.25
$ws.Range("J103").Value = 20593.334
$ws.Range("L103").Value = 20593.334
$ws.Range("N103").Value = -22937.334
